$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5306265354156494
$ws.Range("B1").Value = 1.608950734138489
$ws.Range("C1").Value = 2.327079057693481
$ws.Range("D1").Value = 1.971695423126221
$ws.Range("E1").Value = 0.9562870264053345
